# "Refresh as of the next day" update.
#
# Column D = total days (总天), E = remaining days (剩余),
# F = cycle start date as an 8-digit YYYYMMDD integer (开始时间).
#
# The sheet encodes, for every row: E = D - (today - F) i.e. remaining
# days counts down from D as "today" advances past the start date F.
# This edit advances "today" from 2025-11-08 to 2025-11-09, so every
# row's E drops by 1 day of elapsed time. When a row's countdown would
# reach 0 (E was 1), a new cycle starts instead: F jumps forward by 10
# days and E resets to 10.

function DaysInMonth($y, $m) {
    $dim = @(31,28,31,30,31,30,31,31,30,31,30,31)
    if ((($y % 4) -eq 0 -and ($y % 100) -ne 0) -or (($y % 400) -eq 0)) {
        $dim[1] = 29
    }
    return $dim[$m - 1]
}

# Ordinal day count (within-era, monotonic) for y-m-d, used only for
# differences between two dates in the same neighbourhood.
function ToOrdinal($y, $m, $d) {
    $dim = @(31,28,31,30,31,30,31,31,30,31,30,31)
    if ((($y % 4) -eq 0 -and ($y % 100) -ne 0) -or (($y % 400) -eq 0)) {
        $dim[1] = 29
    }
    $total = $d
    for ($i = 0; $i -lt ($m - 1); $i++) {
        $total = $total + $dim[$i]
    }
    $total = $total + ($y * 366)
    return $total
}

# Add $n days (n can be negative) to y-m-d, return array(y,m,d).
function AddDaysYMD($y, $m, $d, $n) {
    $d = $d + $n
    while ($d -gt (DaysInMonth $y $m)) {
        $d = $d - (DaysInMonth $y $m)
        $m = $m + 1
        if ($m -gt 12) {
            $m = 1
            $y = $y + 1
        }
    }
    while ($d -lt 1) {
        $m = $m - 1
        if ($m -lt 1) {
            $m = 12
            $y = $y - 1
        }
        $d = $d + (DaysInMonth $y $m)
    }
    $result = @($y, $m, $d)
    return $result
}

function YmdToInt($y, $m, $d) {
    $s = "{0}{1:D2}{2:D2}" -f $y, $m, $d
    return [int]$s
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Today" reference used by the sheet: before this edit it was
# 2025-11-08 (consistent with every row's D/E/F), this edit moves it
# one day forward.
$todayY = 2025
$todayM = 11
$todayD = 9
$todayOrd = ToOrdinal $todayY $todayM $todayD

$lastRow = $ws.UsedRange.Rows.Count
if (-not $lastRow -or $lastRow -lt 2) {
    $lastRow = 99
}

for ($row = 2; $row -le $lastRow; $row++) {
    $dCell = $ws.Cells.Item($row, 4)
    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)

    $dVal = $dCell.Value2
    $fVal = $fCell.Value2

    if ($null -eq $dVal -or $null -eq $fVal) {
        continue
    }

    $fStr = "{0}" -f ([int64]$fVal)
    if ($fStr.Length -ne 8) {
        # Malformed / non-date start value (e.g. data-entry typo) - leave as-is.
        continue
    }

    $fy = [int]$fStr.Substring(0, 4)
    $fm = [int]$fStr.Substring(4, 2)
    $fd = [int]$fStr.Substring(6, 2)

    $fOrd = ToOrdinal $fy $fm $fd
    $elapsed = $todayOrd - $fOrd
    $newE = [int]$dVal - $elapsed

    if ($newE -le 0) {
        # Countdown exhausted: roll the cycle start forward by 10 days
        # and recompute remaining days against the new start.
        $parts = AddDaysYMD $fy $fm $fd 10
        $fy = $parts[0]
        $fm = $parts[1]
        $fd = $parts[2]
        $fOrd = ToOrdinal $fy $fm $fd
        $elapsed = $todayOrd - $fOrd
        $newE = [int]$dVal - $elapsed
    }

    $eCell.Value2 = $newE
    $fCell.Value2 = (YmdToInt $fy $fm $fd)
}
